{"js": "// Update the date heading and refresh the 25 division problems/answers in\n// the 5x5 grid of filled rows (rows 0, 4, 8, 12, 16 of the 20-row table).\n\n// 1) Update the \"YYYY-MM-DD Weekday\" heading paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst headingRange = paragraphs.items[0].getRange();\nheadingRange.insertText(\"2026-01-13 Tuesday\", Word.InsertLocation.replace);\n\n// 2) Update every data cell in the table with the new division problems.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices (0-based, within the 20-row table) that hold visible content;\n// the rows in between are blank spacer rows.\nconst rowIndices = [0, 4, 8, 12, 16];\n\n// New cell text, 5 values per populated row, left-to-right.\nconst newValues = [\n  [\"51\u00f78=6, 3\", \"96\u00f78=12, 0\", \"63\u00f72=31, 1\", \"96\u00f78=12, 0\", \"45\u00f78=5, 5\"],\n  [\"72\u00f79=8, 0\", \"61\u00f78=7, 5\", \"57\u00f72=28, 1\", \"48\u00f79=5, 3\", \"73\u00f75=14, 3\"],\n  [\"92\u00f77=13, 1\", \"46\u00f78=5, 6\", \"84\u00f74=21, 0\", \"81\u00f74=20, 1\", \"37\u00f75=7, 2\"],\n  [\"46\u00f79=5, 1\", \"75\u00f74=18, 3\", \"50\u00f78=6, 2\", \"51\u00f73=17, 0\", \"68\u00f76=11, 2\"],\n  [\"16\u00f79=1, 7\", \"35\u00f73=11, 2\", \"39\u00f79=4, 3\", \"26\u00f76=4, 2\", \"89\u00f78=11, 1\"],\n];\n\nfor (let r = 0; r < rowIndices.length; r++) {\n  const tableRow = rowIndices[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(tableRow, c);\n    const cellRange = cell.body.paragraphs.getFirst().getRange();\n    cellRange.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and refresh the 25 division problems/answers in\n# the 5x5 grid of filled rows (table rows 1, 5, 9, 13, 17 of the 20-row\n# table, using Word's 1-based row/column numbering).\n\n$d = $word.ActiveDocument\n\n# 1) Update the \"YYYY-MM-DD Weekday\" heading paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-13 Tuesday\"\n\n# 2) Update every data cell in the table with the new division problems.\n$t = $d.Tables.Item(1)\n\n# Row indices (1-based) that hold visible content; the rows in between are\n# blank spacer rows.\n$rowIndices = @(1, 5, 9, 13, 17)\n\n# New cell text, 5 values per populated row, left-to-right.\n$newValues = @(\n    @(\"51\u00f78=6, 3\", \"96\u00f78=12, 0\", \"63\u00f72=31, 1\", \"96\u00f78=12, 0\", \"45\u00f78=5, 5\"),\n    @(\"72\u00f79=8, 0\", \"61\u00f78=7, 5\", \"57\u00f72=28, 1\", \"48\u00f79=5, 3\", \"73\u00f75=14, 3\"),\n    @(\"92\u00f77=13, 1\", \"46\u00f78=5, 6\", \"84\u00f74=21, 0\", \"81\u00f74=20, 1\", \"37\u00f75=7, 2\"),\n    @(\"46\u00f79=5, 1\", \"75\u00f74=18, 3\", \"50\u00f78=6, 2\", \"51\u00f73=17, 0\", \"68\u00f76=11, 2\"),\n    @(\"16\u00f79=1, 7\", \"35\u00f73=11, 2\", \"39\u00f79=4, 3\", \"26\u00f76=4, 2\", \"89\u00f78=11, 1\")\n)\n\nfor ($r = 0; $r -lt $rowIndices.Length; $r++) {\n    $tableRow = $rowIndices[$r]\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $t.Cell($tableRow, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n"}
